# 5.5.1 data sheet: add a new "2021" column (R) mirroring the existing
# "2020" column (Q), then move the active selection as recorded by the
# author when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for column R (year header + the data point).
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 20.5

# Column R should look exactly like column Q (same borders / font /
# number format / alignment) - copy Q4:Q5's formatting onto R4:R5
# without disturbing the values we just set.
$ws.Range("Q4:Q5").Copy()
$ws.Range("R4:R5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the single-cell selection that was recorded in the saved file.
$ws.Range("S12").Select()
